# Adds "Top 5 Sector" / "Top 5 %NAV" columns (U:V) to the ETF extract sheet,
# matching the header formatting used by the rest of row 1, then widens the
# affected rows so the new wrapped text is fully visible.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (U1:V1) ---------------------------------------------------
# Copy the look of the existing header cells (bold white-on-blue, bordered,
# wrap text) onto U1:V1 before writing the new header text.
$ws.Range("A1:B1").Copy()
$ws.Range("U1:V1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("U1").Value = "Top 5 Sector"
$ws.Range("V1").Value = "Top 5 %NAV"

# --- Data rows (U2:V6) --------------------------------------------------
$ws.Range("U2").Value = "Technology,Health Care,Industrials,Consumer Staples,Financials"
$ws.Range("V2").Value = "17.7,12.0,12.0,11.1,7.1"

$ws.Range("U3").Value = "Financials,Consumer Discretionary,Government,Industrials,Technology"
$ws.Range("V3").Value = "31.5,5.3,2.3,1.79,1.62"

$ws.Range("U4").Value = "Financials,Health Care,Consumer Staples,Consumer Discretionary,Technology"
$ws.Range("V4").Value = "19.1,12.7,11.2,10.6,9.9"

$ws.Range("U5").Value = "Technology,Communications,Consumer Discretionary,Health Care,Consumer Staples"
$ws.Range("V5").Value = "41.5,14.4,10.8,5.8,5.4"

$ws.Range("U6").Value = "Technology,Health Care,Consumer Discretionary,Communications,Consumer Staples"
$ws.Range("V6").Value = "35.2,13.4,10.4,7.9,6.3"

# --- Row heights grow to fit the newly-wrapped sector/NAV text ---------
$ws.Rows.Item(2).RowHeight = 93.6
$ws.Rows.Item(3).RowHeight = 93.6
$ws.Rows.Item(4).RowHeight = 120
$ws.Rows.Item(5).RowHeight = 120
$ws.Rows.Item(6).RowHeight = 120

# --- Selection moves to the newly added column -------------------------
$ws.Range("U7").Select()
